$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'23.110.68"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -3.16%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.601.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -3.16%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.13%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'  -0.10%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'301.65"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -3.06%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.3784"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -2.83%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3648"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -4.43%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'49.87"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -3.25%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  -6.33%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.08167"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.87%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.15%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'23.02"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -4.38%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.595"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -6.63%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.00001261"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -4.28%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'7.401"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -8.61%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.599.76"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -3.30%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'91.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -2.55%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06854"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.26%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'18.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -7.33%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'6.594"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -5.62%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.5557"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -6.52%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -5.27%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'23.104.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -3.16%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.343"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -3.74%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'2.739"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -7.67%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'21.15"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -4.33%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'149.84"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.75%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'5.263"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -3.30%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'132.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.98%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.392"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Value = "'6.845"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -13.22%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'1.777.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -3.13%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.9553"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.94%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.07695"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -6.27%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'6.276"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.65%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.2556"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -4.73%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.02720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -6.51%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.08896"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -2.88%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'10.07"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -6.88%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.371"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.68%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.7084"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -6.75%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'12.61"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -7.24%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'15.38"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -7.02%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.6613"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -4.85%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'2.319"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.62%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  -0.14%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'3.997"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -2.48%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'131.37"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.16%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.07941"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -4.44%  "
$ws.Range("E51").Style = "Normal"
